# Update 24-Mei-2021, end of day update.
#
# The petty-cash "Sheet1" ledger is rolled forward to 24-May-2021:
#  - the opening balance (E2 / "SALDO AWAL") becomes the prior day's closing
#    balance, 1,124,025
#  - the first dated row (row 3) now records 24-May-2021 (serial 44340)
#    instead of 17-May-2021 (serial 44333)
#  - all of the now-stale transaction rows for the days in between
#    (rows 4-32: dates, descriptions, debit/credit amounts) are cleared out,
#    leaving only the running-balance formulas in column E, which recompute
#    to the new carried-forward balance automatically

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New opening balance carried in from the previous day's closing balance.
$ws.Range("E2").Value = 1124025

# Roll the ledger's first transaction date forward to 24-May-2021.
$ws.Range("A3").Value = 44340

# Clear out the now-obsolete transaction detail (dates/descriptions/amounts)
# for the days that have since been consolidated; the shared E-column
# formulas remain and recalculate the running balance.
$ws.Range("A4:D32").Clear()

# End-of-day cursor: frozen pane now starts at row 3, selection on B4.
[void]$ws.Range("B4").Select()
